$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) The run holding the single space between "{{date_time}}" and "{{goal}}"
#    loses its explicit en-US language mark (it reverts to the document's
#    default language). Isolate just that one space character and retype it
#    so it picks up the neighbouring (language-less) formatting instead of
#    carrying over the old explicit <w:lang val="en-US"/>.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("time}} {{goal") | Out-Null
$spaceStart = $anchor.Start + 6
$spaceEnd = $spaceStart + 1
$spaceRange = $d.Range($spaceStart, $spaceEnd)
$spaceRange.Delete()
$insPoint = $d.Range($spaceStart, $spaceStart)
$insPoint.InsertAfter(" ")

# ---------------------------------------------------------------------------
# 2) Extend / reword the "chistota i poryadok" sentence: singular "sobraniya"
#    instead of plural "sobranij", plus the new clause about individual
#    protection measures (up to 50 people).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "собраний.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "собрания, а также гарантируем соблюдение мер индивидуальной защиты (до 50 человек).",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Ответственный" becomes "Ответственный(-ая)" -- insert the feminine
#    suffix marker right after the existing "ый" run.
# ---------------------------------------------------------------------------
$resp = $d.Content
$resp.Find.Execute("Ответственный") | Out-Null
$resp.Collapse(0)
$resp.InsertAfter("(-ая)")
